# Add a new player row (Jordan Clarkson), change T.J. McConnell -> Kevin Porter Jr.
# (with his new team), and move Joel Embiid down below Shai Gilgeous-Alexander.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired table contents (rows 2..18), in order.
$data = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Kevin Porter Jr.", "PG", "LA Clippers"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Dennis Schröder", "PG", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
